# Update the QC debrief/errors sheet with the latest participant data.
# Old rows (sub_001..sub_006, all FALSE, no reason) are replaced with the
# new set of flagged participants and their reasons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (keep header in row 1) before writing the
# smaller replacement table.
$ws.Range("A2:C7").ClearContents()

# Header row stays the same text, just re-assert it.
$ws.Range("A1").Value = "participant"
$ws.Range("B1").Value = "qc_fail_manual"
$ws.Range("C1").Value = "reason"

# New data rows.
$ws.Range("A2").Value = "sub_012"
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = "technical"

$ws.Range("A3").Value = "sub_014"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "technical"

$ws.Range("A4").Value = "sub_020"
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = "missing input for listing names of pictures for the last page."

# Match the saved selection/active cell from the authored workbook.
$ws.Range("C5").Select() | Out-Null
